$d = $word.ActiveDocument

# Mapping of old division expressions to new ones, as per the diff.
$replacements = @(
    @("73÷6=", "63÷4="),
    @("94÷7=", "12÷4="),
    @("68÷6=", "11÷4="),
    @("45÷8=", "41÷5="),
    @("33÷6=", "30÷6="),
    @("87÷6=", "96÷9="),
    @("91÷3=", "71÷3="),
    @("81÷4=", "52÷7="),
    @("88÷8=", "39÷5="),
    @("96÷7=", "66÷8="),
    @("59÷9=", "19÷7="),
    @("88÷2=", "36÷8="),
    @("33÷9=", "50÷7="),
    @("54÷2=", "15÷9="),
    @("34÷4=", "99÷8="),
    @("35÷6=", "47÷3="),
    @("63÷3=", "97÷8="),
    @("96÷6=", "63÷9="),
    @("35÷9=", "95÷3="),
    @("84÷3=", "39÷5="),
    @("54÷6=", "46÷9="),
    @("61÷8=", "19÷3="),
    @("89÷4=", "78÷7="),
    @("69÷6=", "38÷9="),
    @("93÷3=", "51÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
